# Auto-generated Excel COM-interop script to apply market-data refresh diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H, I, J, K, L, M, N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 439
$ws.Range("I6").Value = 298.75
$ws.Range("K6").Value = 896.25
$ws.Range("M6").Value = -784.25
$ws.Range("H29").Value = 2700.111
$ws.Range("H38").Value = 182.83333
$ws.Range("I38").Value = 182.83333
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 548.49999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -176.49999
$ws.Range("N38").ClearContents()
$ws.Range("H43").Value = 1101
$ws.Range("I43").Value = 430
$ws.Range("K43").Value = 430
$ws.Range("M43").Value = -361
$ws.Range("H58").Value = 1820.75
$ws.Range("I58").Value = 1303
$ws.Range("J58").Value = 2056.0908
$ws.Range("K58").Value = 3909
$ws.Range("L58").Value = 6168.2724
$ws.Range("M58").Value = -3759
$ws.Range("N58").Value = -6468.2724
$ws.Range("H76").Value = 4116.1934
$ws.Range("I76").Value = 3721.8262
$ws.Range("J76").Value = 5250
$ws.Range("K76").Value = 3721.8262
$ws.Range("L76").Value = 5250
$ws.Range("M76").Value = -3406.8262
$ws.Range("N76").Value = -5880
$ws.Range("H79").Value = 4116.1934
$ws.Range("I79").Value = 3721.8262
$ws.Range("J79").Value = 5250
$ws.Range("K79").Value = 3721.8262
$ws.Range("L79").Value = 5250
$ws.Range("M79").Value = -2629.8262
$ws.Range("N79").Value = -7434
$ws.Range("H87").Value = 25383.6
$ws.Range("J87").Value = 25383.6
$ws.Range("L87").Value = 25383.6
$ws.Range("N87").Value = -27879.6
$ws.Range("H90").Value = 25383.6
$ws.Range("J90").Value = 25383.6
$ws.Range("L90").Value = 76150.79999999999
$ws.Range("N90").Value = -88630.79999999999
$ws.Range("H100").Value = 1568.762
$ws.Range("I100").Value = 1701.6666
$ws.Range("K100").Value = 1701.6666
$ws.Range("M100").Value = -1160.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19927.93
$ws.Range("I32").Value = 21860.156
$ws.Range("J32").Value = 3504
$ws.Range("K32").Value = 21860.156
$ws.Range("L32").Value = 3504
$ws.Range("M32").Value = -21573.156
$ws.Range("N32").Value = -4078
$ws.Range("H97").Value = 934.43475
$ws.Range("I97").Value = 890.0952
$ws.Range("K97").Value = 890.0952
$ws.Range("M97").Value = -394.0952
$ws.Range("H102").Value = 3258.077
$ws.Range("I102").Value = 2915.5
$ws.Range("K102").Value = 2915.5
$ws.Range("M102").Value = -1293.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1126.2858
$ws.Range("I94").Value = 998
$ws.Range("J94").Value = 1596.6666
$ws.Range("K94").Value = 998
$ws.Range("L94").Value = 1596.6666
$ws.Range("M94").Value = -547
$ws.Range("N94").Value = -2498.6666
$ws.Range("H105").Value = 3941.0967
$ws.Range("I105").Value = 3807.375
$ws.Range("J105").Value = 4399.5713
$ws.Range("K105").Value = 3807.375
$ws.Range("L105").Value = 4399.5713
$ws.Range("M105").Value = -2060.375
$ws.Range("N105").Value = -7893.5713

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 602195.9399999999
$ws.Range("I31").Value = 8874.48
$ws.Range("J31").Value = 1003088.8
$ws.Range("K31").Value = 8874.48
$ws.Range("L31").Value = 1003088.8
$ws.Range("M31").Value = -8579.48
$ws.Range("N31").Value = -1003678.8
$ws.Range("H34").Value = 602195.9399999999
$ws.Range("I34").Value = 8874.48
$ws.Range("J34").Value = 1003088.8
$ws.Range("K34").Value = 8874.48
$ws.Range("L34").Value = 1003088.8
$ws.Range("M34").Value = -8672.48
$ws.Range("N34").Value = -1003492.8
$ws.Range("H99").Value = 1265.0741
$ws.Range("I99").Value = 912.7895
$ws.Range("J99").Value = 2101.75
$ws.Range("K99").Value = 912.7895
$ws.Range("L99").Value = 2101.75
$ws.Range("M99").Value = 585.2105
$ws.Range("N99").Value = -5097.75
$ws.Range("H126").Value = 1265.0741
$ws.Range("I126").Value = 912.7895
$ws.Range("J126").Value = 2101.75
$ws.Range("K126").Value = 2738.3685
$ws.Range("L126").Value = 6305.25
$ws.Range("M126").Value = -268.3685
$ws.Range("N126").Value = -11245.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3376.39
$ws.Range("I68").Value = 1560.35
$ws.Range("J68").Value = 7199.6313
$ws.Range("K68").Value = 4681.049999999999
$ws.Range("L68").Value = 21598.8939
$ws.Range("M68").Value = -3870.049999999999
$ws.Range("N68").Value = -23220.8939
$ws.Range("H71").Value = 3376.39
$ws.Range("I71").Value = 1560.35
$ws.Range("J71").Value = 7199.6313
$ws.Range("K71").Value = 14043.15
$ws.Range("L71").Value = 64796.6817
$ws.Range("M71").Value = -9987.15
$ws.Range("N71").Value = -72908.6817
$ws.Range("H112").Value = 2409
$ws.Range("I112").Value = 1613.5
$ws.Range("J112").Value = 4000
$ws.Range("K112").Value = 4840.5
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = -3732.5
$ws.Range("N112").Value = -14216
$ws.Range("H131").Value = 1131.4386
$ws.Range("I131").Value = 1533
$ws.Range("J131").Value = 1012.7955
$ws.Range("K131").Value = 4599
$ws.Range("L131").Value = 3038.3865
$ws.Range("M131").Value = 441
$ws.Range("N131").Value = -13118.3865

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 4169.3335
$ws.Range("J9").Value = 3754
$ws.Range("L9").Value = 3754
$ws.Range("N9").Value = -4094
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H97").Value = 1422.8
$ws.Range("I97").Value = 1616.3
$ws.Range("J97").Value = 1035.8
$ws.Range("K97").Value = 1616.3
$ws.Range("L97").Value = 1035.8
$ws.Range("M97").Value = -1120.3
$ws.Range("N97").Value = -2027.8
$ws.Range("H123").Value = 10333.111
$ws.Range("J123").Value = 10333.111
$ws.Range("L123").Value = 10333.111
$ws.Range("N123").Value = -15233.111

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2905.0588
$ws.Range("I7").Value = 3291.2222
$ws.Range("K7").Value = 3291.2222
$ws.Range("M7").Value = -3179.2222
$ws.Range("H46").Value = 1029.4375
$ws.Range("I46").Value = 737.9
$ws.Range("J46").Value = 1515.3334
$ws.Range("K46").Value = 737.9
$ws.Range("L46").Value = 1515.3334
$ws.Range("M46").Value = -549.9
$ws.Range("N46").Value = -1891.3334
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H93").Value = 1176.1111
$ws.Range("I93").Value = 1283.5714
$ws.Range("K93").Value = 1283.5714
$ws.Range("M93").Value = -35.57140000000004
$ws.Range("H126").Value = 2905.0588
$ws.Range("I126").Value = 3291.2222
$ws.Range("K126").Value = 9873.6666
$ws.Range("M126").Value = -7403.6666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 29874.5
$ws.Range("I63").Value = 19500
$ws.Range("J63").Value = 40249
$ws.Range("K63").Value = 19500
$ws.Range("L63").Value = 40249
$ws.Range("M63").Value = -18876
$ws.Range("N63").Value = -41497
$ws.Range("H66").Value = 29874.5
$ws.Range("I66").Value = 19500
$ws.Range("J66").Value = 40249
$ws.Range("K66").Value = 58500
$ws.Range("L66").Value = 120747
$ws.Range("M66").Value = -55380
$ws.Range("N66").Value = -126987
